$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force text number-format so the numeric-looking strings are kept as text
# (matching the shared-string / text storage used in the source file)
# rather than being auto-converted to numbers by Excel.
$cells = @("B11", "C11", "D11", "B13", "C13", "D13")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 11: "Enterprises density (per 1000 people)"
$ws.Range("B11").Value = "30.45"
$ws.Range("C11").Value = "0.91"
$ws.Range("D11").Value = "31.36"

# Row 13: "Enterprises (% of total)"
$ws.Range("B13").Value = "97.02"
$ws.Range("C13").Value = "2.89"
$ws.Range("D13").Value = "99.91"
